# Insert two new data rows at row 213 (pushing the existing rows 213.. down
# by two, so the former row 213 becomes row 215, etc., and the last row
# (old 288) becomes row 290), then populate the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(213).EntireRow.Insert()
$ws.Rows.Item(213).EntireRow.Insert()

# New row 213
$ws.Range("A213").Value = 4
$ws.Range("B213").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value = "Los Lagos"
$ws.Range("D213").Value = 44524
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = "Fruta"
$ws.Range("G213").Value = 100108
$ws.Range("H213").Value = "Tropicales y subtropicales"
$ws.Range("I213").Value = 100108006
$ws.Range("J213").Value = "Plátano"
$ws.Range("K213").Value = "Barraganete"
$ws.Range("L213").Value = "Primera"
$ws.Range("M213").Value = 80
$ws.Range("N213").Value = 32000
$ws.Range("O213").Value = 33000
$ws.Range("P213").Value = 32500
$ws.Range("Q213").Value = "$/caja 20 kilos"
$ws.Range("R213").Value = "Ecuador"
$ws.Range("S213").Value = 1625
$ws.Range("T213").Value = 20

# New row 214
$ws.Range("A214").Value = 4
$ws.Range("B214").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C214").Value = "Los Lagos"
$ws.Range("D214").Value = 44524
$ws.Range("E214").Value = 10
$ws.Range("F214").Value = "Fruta"
$ws.Range("G214").Value = 100108
$ws.Range("H214").Value = "Tropicales y subtropicales"
$ws.Range("I214").Value = 100108006
$ws.Range("J214").Value = "Plátano"
$ws.Range("K214").Value = "Sin especificar"
$ws.Range("L214").Value = "Primera Pintón"
$ws.Range("M214").Value = 600
$ws.Range("N214").Value = 23000
$ws.Range("O214").Value = 24000
$ws.Range("P214").Value = 23500
$ws.Range("Q214").Value = "$/caja 20 kilos"
$ws.Range("R214").Value = "Ecuador"
$ws.Range("S214").Value = 1175
$ws.Range("T214").Value = 20
